$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 74-82 (9 trailing rows removed)
$ws.Rows("74:82").Delete()

# Update column B (value) for rows 2-73 with the new figures
$ws.Range("B2").Value = 0.25
$ws.Range("B3").Value = 0.25
$ws.Range("B4").Value = 0.1500000000000057
$ws.Range("B5").Value = 0.2000000000000028
$ws.Range("B6").Value = 0.4000000000000057
$ws.Range("B7").Value = 0.4999999999999858
$ws.Range("B8").Value = 0.09999999999999432
$ws.Range("B9").Value = 0.4000000000000057
$ws.Range("B10").Value = -0.5
$ws.Range("B11").Value = 0.4999999999999858
$ws.Range("B12").Value = 0.2999999999999829
$ws.Range("B13").Value = 0.2250000000000085
$ws.Range("B14").Value = 0.4
$ws.Range("B15").Value = 0.2
$ws.Range("B16").Value = 0.2
$ws.Range("B17").Value = 0.3
$ws.Range("B18").Value = 0.3
$ws.Range("B19").Value = 0.3
$ws.Range("B20").Value = 0.4
$ws.Range("B21").Value = 0.3
$ws.Range("B22").Value = 0.1
$ws.Range("B23").Value = 0.2
$ws.Range("B24").Value = 0.2
$ws.Range("B25").Value = 0.3
$ws.Range("B26").Value = 0.3
$ws.Range("B27").Value = 0.4
$ws.Range("B28").Value = 0.3490000000000038
$ws.Range("B29").Value = 0.4399999999999977
$ws.Range("B30").Value = 0.4399999999999977
$ws.Range("B31").Value = 0.3500000000000085
$ws.Range("B32").Value = 0.4499999999999886
$ws.Range("B33").Value = 0.4000000000000057
$ws.Range("B34").Value = 0.4202440737484352
$ws.Range("B35").Value = 0.4753723183093825
$ws.Range("B36").Value = 0.2999999999999687
$ws.Range("B37").Value = 0.4499999999999886
$ws.Range("B38").Value = 0.4000000000000057
$ws.Range("B39").Value = 0.3500000000000085
$ws.Range("B40").Value = 0.4499999999999744
$ws.Range("B41").Value = 0.4250000000000114
$ws.Range("B42").Value = 0.4000000000000057
$ws.Range("B43").Value = 0.4000000000000057
$ws.Range("B44").Value = 0.4000000000000057
$ws.Range("B45").Value = 0.4000000000000057
$ws.Range("B46").Value = 0.3
$ws.Range("B47").Value = 0.539999999999992
$ws.Range("B48").Value = -0.3
$ws.Range("B49").Value = 0.09999999999996589
$ws.Range("B50").Value = 0.1399999999999864
$ws.Range("B51").Value = -0.09999999999996589
$ws.Range("B52").Value = -10.8
$ws.Range("B53").Value = 7.25
$ws.Range("B54").Value = -0.7999999999999972
$ws.Range("B55").Value = -2.400000000000006
$ws.Range("B56").Value = 3.480874220397794
$ws.Range("B57").Value = 6.400000000000006
$ws.Range("B58").Value = -1.099999999999994
$ws.Range("B59").Value = 0.4209467346675666
$ws.Range("B60").Value = 1.540000000000006
$ws.Range("B61").Value = -0.1490000000000009
$ws.Range("B62").Value = -0.4999966213670604
$ws.Range("B63").Value = -0.7399999999999949
$ws.Range("B64").Value = -0.2510000000000048
$ws.Range("B65").Value = 0.2000000000000028
$ws.Range("B66").Value = 0.1200000000000045
$ws.Range("B67").Value = 0.09999999999999432
$ws.Range("B68").Value = 0
$ws.Range("B69").Value = 0.09999999999999432
$ws.Range("B70").Value = 0.3200015876295765
$ws.Range("B71").Value = 0.09999771493470178
$ws.Range("B72").Value = 0.2000034419242951
$ws.Range("B73").Value = 0.03999999999999204
